$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, matching the source data
# (many of the updated prices would otherwise be auto-parsed as numbers by Excel,
# e.g. "1.00" -> 1, "135.00" -> 135, losing the original text formatting).
$ws.Range("D2:D51").NumberFormat = "@"

# Apply the updated cryptocurrency data
$ws.Range("D2").Value = "60.373.04"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "2.605.04"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "578.31"
$ws.Range("E5").Value = "  +3.78%  "
$ws.Range("D6").Value = "143.20"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "2.612.16"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "6.57"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "0.156"
$ws.Range("E12").Value = "  -3.10%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "3.060.14"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "24.41"
$ws.Range("E15").Value = "  +4.14%  "
$ws.Range("D16").Value = "60.361.83"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "2.606.92"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "11.46"
$ws.Range("E19").Value = "  +8.11%  "
$ws.Range("D20").Value = "4.63"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "347.32"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").Value = "6.90"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "0.529"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "63.18"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").Value = "1.87"
$ws.Range("E30").Value = "  +10.16%  "
$ws.Range("D31").Value = "6.39"
$ws.Range("E31").Value = "  +3.66%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "166.17"
$ws.Range("E33").Value = "  +4.78%  "
$ws.Range("D34").Value = "19.45"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "4.29"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("D36").Value = "1.30"
$ws.Range("E36").Value = "  +9.92%  "
$ws.Range("D37").Value = "0.990"
$ws.Range("E37").Value = "  +7.86%  "
$ws.Range("E38").Value = "  +7.41%  "
$ws.Range("D39").Value = "38.08"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").Value = "314.71"
$ws.Range("E40").Value = "  +7.91%  "
$ws.Range("E41").Value = "  +5.50%  "
$ws.Range("D42").Value = "0.839"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "135.00"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "0.0996"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "19.86"
$ws.Range("E46").Value = "  +3.61%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0553"
$ws.Range("E47").Value = "  +2.54%  "
$ws.Range("D48").Value = "5.02"
$ws.Range("E48").Value = "  +9.42%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.606"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "20.04"
$ws.Range("E50").Value = "  +5.33%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0241"
$ws.Range("E51").Value = "  +0.68%  "
